# Insert a new bullet "Add initials of data person in additional column"
# right after the "Create unique trial/individual ID" bullet (same list,
# same indent level) and before "Filter out trials with salinity issue ...".

$d = $word.ActiveDocument

# Locate the anchor paragraph via Find so we don't depend on a hard-coded
# paragraph index.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Create unique trial/individual ID",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate anchor paragraph 'Create unique trial/individual ID'"
}

$anchorIndex = $findRange.Paragraphs.Item(1).Index
$anchorPara = $d.Paragraphs.Item($anchorIndex)

# Insert a new paragraph after the anchor; it inherits the anchor's
# paragraph formatting (ListParagraph style, numPr ilvl 1 / numId 6, and
# the Roboto/20/20/en-US run formatting), matching the target markup.
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Add initials of data person in additional column"
